$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "245.85"
Set-TextValue "G2" "4"
Set-TextValue "D3" "25.37"
Set-TextValue "G3" "4"
Set-TextValue "D4" "5.144"
Set-TextValue "G4" "4"
Set-TextValue "D5" "0.05573"
Set-TextValue "G5" "4"
Set-TextValue "D6" "6.509"
Set-TextValue "G6" "4"
Set-TextValue "D7" "3.019"
Set-TextValue "G7" "4"
Set-TextValue "D8" "0.8186"
Set-TextValue "G8" "4"
Set-TextValue "D9" "0.8493"
Set-TextValue "G9" "4"
Set-TextValue "D10" "0.1342"
Set-TextValue "G10" "4"
Set-TextValue "D11" "0.06960"
Set-TextValue "G11" "4"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.02882"
Set-TextValue "E12" "11BitrueCoinBTR"
Set-TextValue "G12" "4"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09381"
Set-TextValue "E13" "12BitMartTokenBMX"
Set-TextValue "G13" "4"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001516"
Set-TextValue "E14" "13BitForexTokenBF"
Set-TextValue "G14" "4"
Set-TextValue "B15" "One"
Set-TextValue "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0005963"
Set-TextValue "E15" "14OneONE"
Set-TextValue "G15" "4"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.006100"
Set-TextValue "E16" "15TigerCashTCH"
Set-TextValue "G16" "4"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.501"
Set-TextValue "E17" "16LEOLEO"
Set-TextValue "G17" "4"
Set-TextValue "B18" "BTSEToken"
Set-TextValue "C18" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D18" "2.063"
Set-TextValue "E18" "17BTSETokenBTSE"
Set-TextValue "G18" "4"
Set-TextValue "B19" "BitpandaEcosystemToken"
Set-TextValue "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D19" "0.3179"
Set-TextValue "E19" "18BitpandaEcosystemTokenBEST"
Set-TextValue "G19" "4"
Set-TextValue "B20" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C20" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D20" "0.03165"
Set-TextValue "E20" "19LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "G20" "4"
Set-TextValue "G21" "4"
Set-TextValue "D22" "3.759"
Set-TextValue "G22" "4"
Set-TextValue "D23" "0.04731"
Set-TextValue "G23" "4"
Set-TextValue "G24" "4"
Set-TextValue "D25" "0.001250"
Set-TextValue "G25" "4"
Set-TextValue "D26" "0.004640"
Set-TextValue "G26" "4"
Set-TextValue "D27" "0.00009703"
Set-TextValue "G27" "4"
Set-TextValue "E28" "27UpBotsUBXTWorstin24h"
Set-TextValue "G28" "4"
Set-TextValue "G29" "4"
Set-TextValue "G30" "4"
Set-TextValue "G31" "4"
Set-TextValue "G32" "4"
Set-TextValue "G33" "4"
Set-TextValue "G34" "4"
Set-TextValue "G35" "4"
Set-TextValue "G36" "4"
Set-TextValue "G37" "4"
Set-TextValue "G38" "4"
Set-TextValue "G39" "4"
Set-TextValue "G40" "4"
Set-TextValue "D41" "0.1363"
Set-TextValue "G41" "4"
Set-TextValue "B42" "KickToken"
Set-TextValue "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.006192"
Set-TextValue "E42" "41KickTokenKICK"
Set-TextValue "G42" "4"
Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002625"
Set-TextValue "E43" "42CEJICEJI"
Set-TextValue "G43" "4"
Set-TextValue "D44" "0.008301"
Set-TextValue "G44" "4"
Set-TextValue "D45" "0.00005301"
Set-TextValue "G45" "4"
Set-TextValue "G46" "4"
Set-TextValue "D47" "0.1891"
Set-TextValue "G47" "4"
Set-TextValue "D48" "0.002122"
Set-TextValue "G48" "4"
Set-TextValue "G49" "4"
Set-TextValue "G50" "4"
Set-TextValue "G51" "4"
